$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract = @'
Background
A cluster of pneumonia cases were reported by Wuhan Municipal Health Commission, China in December 2019. A novel coronavirus was eventually identified, and became the COVID-19 epidemic that affected public health and life.
 We investigated the psychological status and behavior changes of the general public in China from January 30 to February 3, 2020.
Methods
Respondents were recruited via social media (WeChat) and completed an online questionnaire.
 We used the State-Trait Anxiety Inventory, Self-rating Depression Scale, and Symptom Checklist-90 to evaluate psychological status.
 We also investigated respondents’ behavior changes.
 Quantitative data were analyzed by t-tests or analysis of variance, and classified data were analyzed with chi-square tests.
Results
In total, 608 valid questionnaires were obtained.
 More respondents had state anxiety than trait anxiety (15.8% vs 4.0%).
 Depression was found among 27.1% of respondents and 7.7% had psychological abnormalities.
 About 10.1% of respondents suffered from phobia.
 Our analysis of the relationship between subgroup characteristics and psychological status showed that age, gender, knowledge about COVID-19, degree of worry about epidemiological infection, and confidence about overcoming the outbreak significantly influenced psychological status.
 Around 93.3% of respondents avoided going to public places and almost all respondents reduced Spring Festival-related activities.
 At least 70.9% of respondents chose to take three or more preventive measures to avoid infection.
 The three most commonly used prevention measures were making fewer trips outside and avoiding contact (98.0%), wearing a mask (83.7%), and hand hygiene (82.4%).
Conclusions
We need to pay more attention to public psychological stress, especially among young people, as they are likely to experience anxiety, depression, and psychological abnormalities.
 Different psychological interventions could be formulated according to the psychological characteristics of different gender and age groups.
 The majority of respondents followed specific behaviors required by the authorities, but it will take time to observe the effects of these behaviors on the epidemic.

'@

$authors = @'
[Xi%Liu%NULL%1,       Wen-Tao%Luo%NULL%1,       Ying%Li%NULL%1,       Chun-Na%Li%NULL%1,       Zhong-Si%Hong%NULL%1,       Hui-Li%Chen%NULL%1,       Fei%Xiao%NULL%0,       Jin-Yu%Xia%xiajinyu@mail.sysu.edu.cn%1]
'@

$ws.Range("D6").Value = $abstract
$ws.Range("E6").Value = $authors
